$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two image-filename cells with refreshed data-pull dates ---
# D6: commercial revenue-per-vessel chart (2023 -> 2024 data year, new pull date)
$ws.Range("D6").Value = "AVGVESREVperYr_BLACK_SEABASS_2024_DOLlb_2025-04-17.png"
# D7: number-of-commercial-vessels chart (new pull date)
$ws.Range("D7").Value = "N_Commercial_Vessels_Landing_BLACK_SEABASS_2025-04-17.png"

# D6 picked up a distinct font treatment (monospace, vertically centered)
$f = $ws.Range("D6").Font
$f.Name = "Consolas"
$f.Size = 10
$ws.Range("D6").VerticalAlignment = -4108

# --- Row heights settle slightly after the edit (re-wrap on newer Excel) ---
$ws.Rows(2).RowHeight = 77.25
$ws.Rows(3).RowHeight = 90
$ws.Rows(4).RowHeight = 51.75
$ws.Rows(5).RowHeight = 39
$ws.Rows(6).RowHeight = 39
$ws.Rows(7).RowHeight = 39

# --- Column widths also nudge slightly, and column D gets an explicit width ---
$ws.Columns(1).ColumnWidth = 13.43
$ws.Columns(2).ColumnWidth = 16.26
$ws.Columns(3).ColumnWidth = 50.75
$ws.Columns(4).ColumnWidth = 60.75

# --- View: drop the old scroll position and move the live selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
